$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Telefone"
$ws.Range("I1").Value = "linkedin"
$ws.Range("E1").Value = "E-mail"

$ws.Range("E1").Select()
